$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the overall Pass/Fail summary formula in B6: the SUM range shrinks
#    from H8:H34 to H8:H33 (row 34 now holds a different test after the
#    row re-order performed below).
# ---------------------------------------------------------------------------
$ws.Range("B6").Formula = '=IF(SUM(H8:H33)=23,"Pass","Fail")'

# ---------------------------------------------------------------------------
# 2. Re-order rows 32:34.
#    Row 32 ("amp_Imbalance", the MAX(B28:B31)-MIN(B28:B31) summary row) is
#    moved down so it becomes the new row 34.
#    Row 33 ("Res_TxMPDX") becomes the new row 32.
#    Row 34 ("Res_TxMPDY") becomes the new row 33.
#    The new layout is written explicitly cell-by-cell (instead of using
#    Cut/Insert) so formulas, value types and styles all end up exactly
#    right.
# ---------------------------------------------------------------------------

# --- New row 32 (was row 33: Res_TxMPDX) -----------------------------------
$ws.Range("A32").Value = "Res_TxMPDX"
$ws.Range("B32").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("E32").Formula = '=IF(COUNTA(B32:D32),B32,"NA")'
$ws.Range("F32").Value = 0.00029999999999999997
$ws.Range("G32").Value = 0.0050000000000000001
$ws.Range("H32").Formula = '=IF(I32="PASS",1,0)'
$ws.Range("I32").Formula = '=IF(AND(E32>=F32,E32<G32,NOT( E32 ="NA")),"Pass","Fail")'

# --- New row 33 (was row 34: Res_TxMPDY) -----------------------------------
$ws.Range("A33").Value = "Res_TxMPDY"
$ws.Range("B33").ClearContents()
$ws.Range("C33").ClearContents()
$ws.Range("D33").ClearContents()
$ws.Range("E33").Formula = '=IF(COUNTA(B33:D33),B33,"NA")'
$ws.Range("F33").Value = 0.00029999999999999997
$ws.Range("G33").Value = 0.0050000000000000001
$ws.Range("H33").Formula = '=IF(I33="PASS",1,0)'
$ws.Range("I33").Formula = '=IF(AND(E33>=F33,E33<G33,NOT( E33 ="NA")),"Pass","Fail")'

# --- New row 34 (was row 32: amp_Imbalance) --------------------------------
$ws.Range("A34").Value = "amp_Imbalance"
$ws.Range("B34").Formula = "=MAX(B28:B31)-MIN(B28:B31)"
$ws.Range("C34").ClearContents()
$ws.Range("D34").ClearContents()
$ws.Range("E34").Formula = '=IF(COUNTA(B34:D34),B34,"NA")'
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 6
$ws.Range("H34").Formula = '=IF(I34="PASS",1,0)'
$ws.Range("I34").Formula = '=IF(AND(E34>F34,E34<G34),"Pass","Fail")'

# ---------------------------------------------------------------------------
# 3. Update conditional formatting so the FAIL/PASS highlighting follows the
#    re-ordered rows: the rule that used to watch I28:I32 now watches
#    I28:I31 (plus the relocated I34), and the rule that used to watch
#    I33:I34 now watches I32:I33.
# ---------------------------------------------------------------------------

# I28:I32 -> I28:I31 (the 5th row of that old block, I32, moved to I34)
$i2832 = $ws.Range("I28:I32").FormatConditions
$i2832.Item(1).ModifyAppliesToRange($ws.Range("I28:I31"))

# I33:I34 -> I32:I33 (shift the block up by one row)
$i3334 = $ws.Range("I33:I34").FormatConditions
$i3334.Item(1).ModifyAppliesToRange($ws.Range("I32:I33"))

# The relocated row (new I34, old I32) keeps its own FAIL/PASS highlighting;
# give it the same two text-contains rules used elsewhere on this sheet.
$i34Fail = $ws.Range("I34").FormatConditions.Add(9, 0, $null, $null, "FAIL")
$i34Fail.Interior.Color = 255
$i34Pass = $ws.Range("I34").FormatConditions.Add(9, 0, $null, $null, "PASS")
$i34Pass.Interior.Color = 5296274
